$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 27028412
$ws.Cells.Item(40, 9).Value = 1361.5714
$ws.Cells.Item(40, 10).Value = 43479660
$ws.Cells.Item(40, 11).Value = 1361.5714
$ws.Cells.Item(40, 12).Value = 43479660
$ws.Cells.Item(40, 13).Value = -1186.5714
$ws.Cells.Item(40, 14).Value = -43480010

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(64, 8).Value = 2821.4285
$ws.Cells.Item(64, 9).Value = 2793.75
$ws.Cells.Item(64, 10).Value = 2838.4614
$ws.Cells.Item(64, 11).Value = 2793.75
$ws.Cells.Item(64, 12).Value = 2838.4614
$ws.Cells.Item(64, 13).Value = -2545.75
$ws.Cells.Item(64, 14).Value = -3334.4614

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(67, 8).Value = 2821.4285
$ws.Cells.Item(67, 9).Value = 2793.75
$ws.Cells.Item(67, 10).Value = 2838.4614
$ws.Cells.Item(67, 11).Value = 2793.75
$ws.Cells.Item(67, 12).Value = 2838.4614
$ws.Cells.Item(67, 13).Value = -1935.75
$ws.Cells.Item(67, 14).Value = -4554.4614

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(70, 8).Value = 928.26086
$ws.Cells.Item(70, 9).Value = 903.75
$ws.Cells.Item(70, 10).Value = 941.3333
$ws.Cells.Item(70, 11).Value = 2711.25
$ws.Cells.Item(70, 12).Value = 2823.9999
$ws.Cells.Item(70, 13).Value = -2441.25
$ws.Cells.Item(70, 14).Value = -3363.9999

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(73, 8).Value = 928.26086
$ws.Cells.Item(73, 9).Value = 903.75
$ws.Cells.Item(73, 10).Value = 941.3333
$ws.Cells.Item(73, 11).Value = 2711.25
$ws.Cells.Item(73, 12).Value = 2823.9999
$ws.Cells.Item(73, 13).Value = -1775.25
$ws.Cells.Item(73, 14).Value = -4695.9999

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(74, 8).Value = 4172.032
$ws.Cells.Item(74, 9).Value = 4744.525
$ws.Cells.Item(74, 10).Value = 3131.1365
$ws.Cells.Item(74, 11).Value = 4744.525
$ws.Cells.Item(74, 12).Value = 3131.1365
$ws.Cells.Item(74, 13).Value = -3808.525
$ws.Cells.Item(74, 14).Value = -5003.136500000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(77, 8).Value = 4172.032
$ws.Cells.Item(77, 9).Value = 4744.525
$ws.Cells.Item(77, 10).Value = 3131.1365
$ws.Cells.Item(77, 11).Value = 23722.625
$ws.Cells.Item(77, 12).Value = 15655.6825
$ws.Cells.Item(77, 13).Value = -19042.625
$ws.Cells.Item(77, 14).Value = -25015.6825

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(116, 8).Value = 4407.1353
$ws.Cells.Item(116, 9).Value = 4935
$ws.Cells.Item(116, 10).Value = 4004.9524
$ws.Cells.Item(116, 11).Value = 4935
$ws.Cells.Item(116, 12).Value = 4004.9524
$ws.Cells.Item(116, 13).Value = -1493
$ws.Cells.Item(116, 14).Value = -10888.9524

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(125, 8).Value = 2109
$ws.Cells.Item(125, 9).Value = 1032
$ws.Cells.Item(125, 10).Value = 2206.9092
$ws.Cells.Item(125, 11).Value = 9288
$ws.Cells.Item(125, 12).Value = 19862.1828
$ws.Cells.Item(125, 13).Value = -6828
$ws.Cells.Item(125, 14).Value = -24782.1828

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(127, 8).Value = 624.53845
$ws.Cells.Item(127, 9).Value = 343.25
$ws.Cells.Item(127, 10).Value = 4000
$ws.Cells.Item(127, 11).Value = 1029.75
$ws.Cells.Item(127, 12).Value = 12000
$ws.Cells.Item(127, 13).Value = 3930.25
$ws.Cells.Item(127, 14).Value = -21920

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 1436.9375
$ws.Cells.Item(61, 9).Value = 1432.2069
$ws.Cells.Item(61, 10).Value = 1482.6666
$ws.Cells.Item(61, 11).Value = 1432.2069
$ws.Cells.Item(61, 12).Value = 1482.6666
$ws.Cells.Item(61, 13).Value = -1220.2069
$ws.Cells.Item(61, 14).Value = -1906.6666

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(102, 8).Value = 125001736
$ws.Cells.Item(102, 9).Value = 1811.6666
$ws.Cells.Item(102, 10).Value = 500001500
$ws.Cells.Item(102, 11).Value = 1811.6666
$ws.Cells.Item(102, 12).Value = 500001500
$ws.Cells.Item(102, 13).Value = -189.6666
$ws.Cells.Item(102, 14).Value = -500004744

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(122, 8).Value = 1984.12
$ws.Cells.Item(122, 9).Value = 1880.15
$ws.Cells.Item(122, 11).Value = 5640.450000000001
$ws.Cells.Item(122, 13).Value = -3190.450000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(136, 8).Value = 1436.9375
$ws.Cells.Item(136, 9).Value = 1432.2069
$ws.Cells.Item(136, 10).Value = 1482.6666
$ws.Cells.Item(136, 11).Value = 4296.620699999999
$ws.Cells.Item(136, 12).Value = 4447.9998
$ws.Cells.Item(136, 13).Value = -1746.620699999999
$ws.Cells.Item(136, 14).Value = -9547.9998

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(80, 8).Value = 470.6111
$ws.Cells.Item(80, 9).Value = 167.25
$ws.Cells.Item(80, 10).Value = 557.2857
$ws.Cells.Item(80, 11).Value = 167.25
$ws.Cells.Item(80, 12).Value = 557.2857
$ws.Cells.Item(80, 13).Value = 830.75
$ws.Cells.Item(80, 14).Value = -2553.2857

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(83, 8).Value = 470.6111
$ws.Cells.Item(83, 9).Value = 167.25
$ws.Cells.Item(83, 10).Value = 557.2857
$ws.Cells.Item(83, 11).Value = 836.25
$ws.Cells.Item(83, 12).Value = 2786.4285
$ws.Cells.Item(83, 13).Value = 4155.75
$ws.Cells.Item(83, 14).Value = -12770.4285

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 5901.5884
$ws.Cells.Item(107, 9).Value = 7610.091
$ws.Cells.Item(107, 10).Value = 2769.3333
$ws.Cells.Item(107, 11).Value = 7610.091
$ws.Cells.Item(107, 12).Value = 2769.3333
$ws.Cells.Item(107, 13).Value = -5690.091
$ws.Cells.Item(107, 14).Value = -6609.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1609.025
$ws.Cells.Item(31, 9).Value = 1273.1842
$ws.Cells.Item(31, 10).Value = 7990
$ws.Cells.Item(31, 11).Value = 1273.1842
$ws.Cells.Item(31, 12).Value = 7990
$ws.Cells.Item(31, 13).Value = -978.1841999999999
$ws.Cells.Item(31, 14).Value = -8580

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value = 1609.025
$ws.Cells.Item(34, 9).Value = 1273.1842
$ws.Cells.Item(34, 10).Value = 7990
$ws.Cells.Item(34, 11).Value = 1273.1842
$ws.Cells.Item(34, 12).Value = 7990
$ws.Cells.Item(34, 13).Value = -1071.1842
$ws.Cells.Item(34, 14).Value = -8394

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 2656.2
$ws.Cells.Item(58, 9).Value = 2743.3684
$ws.Cells.Item(58, 11).Value = 2743.3684
$ws.Cells.Item(58, 13).Value = -2540.3684

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(62, 8).Value = 2852.4707
$ws.Cells.Item(62, 9).Value = 2777.2727
$ws.Cells.Item(62, 10).Value = 2990.3333
$ws.Cells.Item(62, 11).Value = 2777.2727
$ws.Cells.Item(62, 12).Value = 2990.3333
$ws.Cells.Item(62, 13).Value = -2153.2727
$ws.Cells.Item(62, 14).Value = -4238.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(65, 8).Value = 2852.4707
$ws.Cells.Item(65, 9).Value = 2777.2727
$ws.Cells.Item(65, 10).Value = 2990.3333
$ws.Cells.Item(65, 11).Value = 13886.3635
$ws.Cells.Item(65, 12).Value = 14951.6665
$ws.Cells.Item(65, 13).Value = -10766.3635
$ws.Cells.Item(65, 14).Value = -21191.6665

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(99, 8).Value = 1562.4
$ws.Cells.Item(99, 9).Value = 1562.4
$ws.Cells.Item(99, 11).Value = 1562.4
$ws.Cells.Item(99, 13).Value = -64.40000000000009

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(122, 8).Value = 2361.818
$ws.Cells.Item(122, 9).Value = 1000
$ws.Cells.Item(122, 10).Value = 4745
$ws.Cells.Item(122, 11).Value = 3000
$ws.Cells.Item(122, 12).Value = 14235
$ws.Cells.Item(122, 13).Value = -550
$ws.Cells.Item(122, 14).Value = -19135

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(126, 8).Value = 1562.4
$ws.Cells.Item(126, 9).Value = 1562.4
$ws.Cells.Item(126, 11).Value = 4687.200000000001
$ws.Cells.Item(126, 13).Value = -2217.200000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(136, 8).Value = 2656.2
$ws.Cells.Item(136, 9).Value = 2743.3684
$ws.Cells.Item(136, 11).Value = 8230.1052
$ws.Cells.Item(136, 13).Value = -5680.1052

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(92, 8).Value = 653.13336
$ws.Cells.Item(92, 9).Value = 368.33334
$ws.Cells.Item(92, 10).Value = 1080.3334
$ws.Cells.Item(92, 11).Value = 1105.00002
$ws.Cells.Item(92, 12).Value = 3241.0002
$ws.Cells.Item(92, 13).Value = 142.9999800000001
$ws.Cells.Item(92, 14).Value = -5737.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(137, 8).Value = 1847.3077
$ws.Cells.Item(137, 9).Value = 1455.909
$ws.Cells.Item(137, 10).Value = 4000
$ws.Cells.Item(137, 11).Value = 4367.727000000001
$ws.Cells.Item(137, 12).Value = 12000
$ws.Cells.Item(137, 13).Value = 732.2729999999992
$ws.Cells.Item(137, 14).Value = -22200

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 2318.8333
$ws.Cells.Item(102, 9).Value = 1720.4166
$ws.Cells.Item(102, 10).Value = 3515.6667
$ws.Cells.Item(102, 11).Value = 1720.4166
$ws.Cells.Item(102, 12).Value = 3515.6667
$ws.Cells.Item(102, 13).Value = -98.41660000000002
$ws.Cells.Item(102, 14).Value = -6759.6667

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 6252458.5
$ws.Cells.Item(7, 9).Value = 12502050
$ws.Cells.Item(7, 11).Value = 12502050
$ws.Cells.Item(7, 13).Value = -12501938

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 2350
$ws.Cells.Item(40, 9).Value = 2500
$ws.Cells.Item(40, 11).Value = 2500
$ws.Cells.Item(40, 13).Value = -2364

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 1887.421
$ws.Cells.Item(61, 9).Value = 1659.75
$ws.Cells.Item(61, 11).Value = 1659.75
$ws.Cells.Item(61, 13).Value = -1457.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(113, 8).Value = 1887.421
$ws.Cells.Item(113, 9).Value = 1659.75
$ws.Cells.Item(113, 11).Value = 1659.75
$ws.Cells.Item(113, 13).Value = 510.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(122, 8).Value = 3916.6667
$ws.Cells.Item(122, 10).Value = 4461.5386
$ws.Cells.Item(122, 12).Value = 13384.6158
$ws.Cells.Item(122, 14).Value = -18284.6158

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(126, 8).Value = 6252458.5
$ws.Cells.Item(126, 9).Value = 12502050
$ws.Cells.Item(126, 11).Value = 37506150
$ws.Cells.Item(126, 13).Value = -37503680

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(140, 8).Value = 52792.57
$ws.Cells.Item(140, 10).Value = 52792.57
$ws.Cells.Item(140, 12).Value = 52792.57
$ws.Cells.Item(140, 14).Value = -63152.57

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 9029.179
$ws.Cells.Item(122, 9).Value = 12656.223
$ws.Cells.Item(122, 10).Value = 2500.5
$ws.Cells.Item(122, 11).Value = 37968.669
$ws.Cells.Item(122, 12).Value = 7501.5
$ws.Cells.Item(122, 13).Value = -35518.669
$ws.Cells.Item(122, 14).Value = -12401.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(126, 8).Value = 2899.7693
$ws.Cells.Item(126, 9).Value = 2572.9333
$ws.Cells.Item(126, 10).Value = 3345.4546
$ws.Cells.Item(126, 11).Value = 7718.7999
$ws.Cells.Item(126, 12).Value = 10036.3638
$ws.Cells.Item(126, 13).Value = -5248.7999
$ws.Cells.Item(126, 14).Value = -14976.3638

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 1605.8667
$ws.Cells.Item(132, 9).Value = 1899.1364
$ws.Cells.Item(132, 10).Value = 799.375
$ws.Cells.Item(132, 11).Value = 5697.4092
$ws.Cells.Item(132, 12).Value = 2398.125
$ws.Cells.Item(132, 13).Value = -3167.4092
$ws.Cells.Item(132, 14).Value = -7458.125
